$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New attendance rows for 9 & 10 Aug 2023 (serials 45147 / 45148), appended
# right after the existing data (row 9 = 8 Aug 2023).
# ---------------------------------------------------------------------------

# Row 10 -> 9 Aug 2023
$ws.Range("A10").Value = 45147
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("B10").Value = "PRESENT"
$ws.Range("C10").Value = "PRESENT"
$ws.Range("D10").Value = "PRESENT"
$ws.Range("E10").Value = "PRESENT"
$ws.Range("F10").Value = "PRESENT"
$ws.Range("G10").Value = "ABSENT"
$ws.Range("H10").Value = "PRESENT"
$ws.Range("I10").Value = "ABSENT"
$ws.Range("J10").Value = "ABSENT"

# Row 11 -> 10 Aug 2023
$ws.Range("A11").Value = 45148
$ws.Range("A11").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("B11").Value = "PRESENT"
$ws.Range("C11").Value = "PRESENT"
$ws.Range("D11").Value = "PRESENT"
$ws.Range("E11").Value = "ABSENT"
$ws.Range("F11").Value = "PRESENT"
$ws.Range("G11").Value = "PRESENT"
$ws.Range("H11").Value = "ABSENT"
$ws.Range("I11").Value = "ABSENT"
$ws.Range("J11").Value = "ABSENT"

# ---------------------------------------------------------------------------
# Explanatory comments (same author bucket/style as the existing "Dell"
# notes) on the newly-marked ABSENT cells.
# ---------------------------------------------------------------------------
$ws.Range("G10").AddComment("Dell:" + [char]10 + "Not informed")
$ws.Range("E11").AddComment("Dell:" + [char]10 + "Not well")
$ws.Range("H11").AddComment("Dell:" + [char]10 + "Not Informed")

# ---------------------------------------------------------------------------
# Extend the "Date" column input-message validation down to the new rows,
# keeping the untouched list validation on B:K intact and in its original
# position.
# ---------------------------------------------------------------------------
$ws.Range("B2:K1048576").Validation.Delete()
$ws.Range("A2:A9").Validation.Delete()

$dateValidation = $ws.Range("A2:A11").Validation
$dateValidation.Add(0, 1, 1, "")
$dateValidation.IgnoreBlank = $false

$ws.Range("B2:K1048576").Validation.Add(3, 1, 1, '"PRESENT, ABSENT"')

# ---------------------------------------------------------------------------
# Match the author's final selection.
# ---------------------------------------------------------------------------
$ws.Range("F12").Select()
